# "more edits to introduction" - the blank first row above the header is
# removed; the header row and every data row below it shift up by one,
# and the trailing blank row follows suit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the blank row 1. Excel renumbers every row below it:
#   old row 2  (header)        -> row 1
#   old rows 3-13 (data)       -> rows 2-12
#   old row 14 (trailing blank)-> row 13
$ws.Rows.Item(1).Delete() | Out-Null

# The data block (now rows 3-12, columns B:E) re-touches its alignment so
# the style engine re-derives the cell format without re-applying an
# explicit (no-op) border, collapsing onto the equivalent borderless style
# already present in the table instead of keeping the redundant one.
$dataRange = $ws.Range("B3:E12")
$dataRange.HorizontalAlignment = -4152   # xlRight
$dataRange.VerticalAlignment = -4108     # xlCenter

# Restore the saved cursor position/selection.
$ws.Range("C15").Select() | Out-Null
